$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 11: D11 formula removed, keep value 37 (G11's shared formula master
# automatically renumbers to si=0 once D's shared-formula group disappears)
$ws.Range("D11").Value = 37

# Rows 12-16: replace previously-shared D formulas with plain values, and
# update the underlying E/F counts that drive the (still-formula) G column.
$ws.Range("D12").Value = 30
$ws.Range("E12").Value = 14
$ws.Range("F12").Value = 2

$ws.Range("D13").Value = 34
$ws.Range("E13").Value = 10
$ws.Range("F13").Value = 2

$ws.Range("D14").Value = 41
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 2

$ws.Range("D15").Value = 34
$ws.Range("E15").Value = 10
$ws.Range("F15").Value = 2

$ws.Range("D16").Value = 41
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 2

# Move the active selection to F23 (matches the saved cursor position).
$ws.Range("F23").Select()
